# Applies the "extent reports" commit changes to keywords.xlsx:
#  1. Fix the text in A4 ("Enter email address." -> "Enter email address")
#  2. Remove the now-redundant "Click on inbox Button" row (row 30),
#     which shifts the trailing "Close Browser" row up to row 30.
#  3. Restore the row-height formatting (15.75pt) that Excel applied to
#     rows 15 and 28 ("Click on Account button" rows) as part of the edit.
#  4. Reposition the sheet view / selection the same way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Correct the typo / trailing period in cell A4.
$ws.Range("A4").Value = "Enter email address"

# 2. Delete the whole row that held "Click on inbox Button" (row 30).
#    This automatically shifts every following row up by one, so the
#    former row 31 ("Close Browser") becomes row 30, and Excel drops the
#    now-unused shared strings for the deleted cells on save.
$ws.Rows.Item(30).Delete()

# 3. Apply the custom row height seen on rows 15 and 28 in the final file.
$ws.Rows.Item(15).RowHeight = 15.75
$ws.Rows.Item(28).RowHeight = 15.75

# 4. Match the final sheet view / selection state.
$ws.Activate()
$ws.Range("A28:XFD28").Select()
$excel.ActiveWindow.ScrollRow = 13
